# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# Mirrors the commit: "Added team record to data" - W/L/T live on the
# same sheet (columns AD:AF) instead of a separate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the rest of row 1 (bold, centered,
# bordered). Copy the formatting from the existing last header cell
# (AC1) onto AD1:AF1, then overwrite with the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-60) gets the same team record values.
for ($row = 2; $row -le 60; $row++) {
    $ws.Cells.Item($row, 30).Value = 92   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 70   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
